# LineasClubPersonalBack.xlsx - add 3 new "Adhesion_Linea_Usuario_*" rows
# (MIX/POS/PRE) to the bottom of the TodasLasLineas sheet, mirroring the
# existing Linea/Caso table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TodasLasLineas")

$newRows = @(
    @{ Linea = "Adhesion_Linea_Usuario_MIX"; Caso = 1162786772 },
    @{ Linea = "Adhesion_Linea_Usuario_POS"; Caso = 1156393564 },
    @{ Linea = "Adhesion_Linea_Usuario_PRE"; Caso = 1162678774 }
)

$startRow = 39

# --- Column A: needs a brand-new font/style (Arial 10, RGB 000000) ---
$aFirst = $ws.Cells.Item($startRow, 1)
$aFirst.Value = $newRows[0].Linea
$aFirst.Font.Color = 0
$aFirst.Font.Size = 10
$aFirst.Font.Name = "Arial"

# Propagate that new style to the rest of column A via copy/paste-format
# (avoids minting yet another style per cell).
$aFirst.Copy()
for ($i = 1; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).PasteSpecial(-4122)
}
for ($i = 1; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newRows[$i].Linea
}

# --- Column B: reuse the existing numeric-cell style already used by B38 ---
$ws.Cells.Item(38, 2).Copy()
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).PasteSpecial(-4122)
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newRows[$i].Caso
}

$excel.CutCopyMode = $false

$ws.Range("E43").Select()
